# Add two new changelog rows (25 and 26) documenting the Microsoft Entra ID
# SSO feature, matching the style/format already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: v0.6.0 Dashboard enhancement entry ---
$ws.Range("A25").Value = 46061
$ws.Range("A25").NumberFormat = "yyyy-mm-dd"
$ws.Range("B25").Value = "0.6.0"
$ws.Range("C25").Value = "Feature"
$ws.Range("D25").Value = "Dashboard enhancement: expanded company dashboard (4 stat cards with links, unified alerts merging training/supervision/appraisals, activity feed from audit log). Redesigned personal dashboard (scheduled vs contracted hours chart with month selector, training pie with category dropdown, clickable charts navigating to filtered pages, supervision next-due date, upcoming shifts panel). Shared libs: shift-colors, audit-messages, activity-feed, upcoming-shifts, formatRelativeTime. Backend: self-service rota fallback in RotaController. URL param filtering for my-training and my-rota. Clean chart re-render keys across all graphs."
$ws.Range("E25").Value = "Claude"

# --- Row 26: v0.7.0 Microsoft Entra ID SSO entry ---
$ws.Range("A26").Value = 46061
$ws.Range("A26").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B26").Value = "0.7.0"
$ws.Range("C26").Value = "Feature"
$ws.Range("D26").Value = "Microsoft Entra ID SSO: MSAL redirect login flow, JIT user provisioning with Employee role, auto-match employees by email, auth_method column (password/microsoft/both), nullable password_hash, Users page auth method controls. Status: core flow working, further testing needed."
$ws.Range("E26").Value = "Claude"
